$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Log_Muestras")

$timestamps = @{
    2  = "2025-11-02T02:04:25.521524"
    3  = "2025-11-02T02:04:25.521524"
    4  = "2025-11-02T02:04:25.521524"
    5  = "2025-11-02T02:04:25.521524"
    6  = "2025-11-02T02:04:25.521524"
    7  = "2025-11-02T02:04:25.521524"
    8  = "2025-11-02T02:04:25.521524"
    9  = "2025-11-02T02:04:25.521524"
    10 = "2025-11-02T02:04:25.522605"
    11 = "2025-11-02T02:04:25.522605"
    12 = "2025-11-02T02:04:25.522605"
    13 = "2025-11-02T02:04:25.522605"
    14 = "2025-11-02T02:04:25.523114"
    15 = "2025-11-02T02:04:25.523161"
    16 = "2025-11-02T02:04:25.523161"
    17 = "2025-11-02T02:04:25.523692"
    18 = "2025-11-02T02:04:25.524689"
    19 = "2025-11-02T02:04:25.525645"
    20 = "2025-11-02T02:04:25.526195"
    21 = "2025-11-02T02:04:25.526195"
    22 = "2025-11-02T02:04:25.526950"
    23 = "2025-11-02T02:04:25.526950"
    24 = "2025-11-02T02:04:25.526950"
    25 = "2025-11-02T02:04:25.527489"
    26 = "2025-11-02T02:04:25.528048"
    27 = "2025-11-02T02:04:25.528666"
    28 = "2025-11-02T02:04:25.528666"
    29 = "2025-11-02T02:04:25.528666"
}

foreach ($row in $timestamps.Keys) {
    $ws.Cells.Item($row, 26).Value = $timestamps[$row]
}
